$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: -> Id 111671188
$ws.Range("A19").Value = 111671188
$ws.Range("B19").Value = 78605
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 6462
$ws.Range("F19").Value = "Stuplav"
$ws.Range("G19").Value = "Nephroma bellum"
$ws.Range("H19").Value = "(Spreng.) Tuck."
$ws.Range("Q19").Value = 558216
$ws.Range("R19").Value = 7067869
$ws.Range("Z19").Value = $null
$ws.Range("AB19").Value = $null

# Row 20: -> Id 111670558
$ws.Range("A20").Value = 111670558
$ws.Range("B20").Value = 96346
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 620
$ws.Range("F20").Value = "Skogsfru"
$ws.Range("G20").Value = "Epipogium aphyllum"
$ws.Range("H20").Value = "Sw."
$ws.Range("Q20").Value = 558134
$ws.Range("R20").Value = 7067979
$ws.Range("Z20").Value = $null
$ws.Range("AB20").Value = $null

# Row 21: -> Id 111670510
$ws.Range("A21").Value = 111670510
$ws.Range("B21").Value = 96346
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 620
$ws.Range("F21").Value = "Skogsfru"
$ws.Range("G21").Value = "Epipogium aphyllum"
$ws.Range("H21").Value = "Sw."
$ws.Range("Q21").Value = 558124
$ws.Range("R21").Value = 7067994
$ws.Range("Z21").Value = $null
$ws.Range("AB21").Value = $null

# Row 22: -> Id 111671179
$ws.Range("A22").Value = 111671179
$ws.Range("B22").Value = 78578
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 6458
$ws.Range("F22").Value = "Lunglav"
$ws.Range("G22").Value = "Lobaria pulmonaria"
$ws.Range("H22").Value = "(L.) Hoffm."
$ws.Range("Q22").Value = 558216
$ws.Range("R22").Value = 7067868
$ws.Range("Z22").Value = $null
$ws.Range("AB22").Value = $null

# Row 23: -> Id 111670477
$ws.Range("A23").Value = 111670477
$ws.Range("B23").Value = 96346
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 620
$ws.Range("F23").Value = "Skogsfru"
$ws.Range("G23").Value = "Epipogium aphyllum"
$ws.Range("H23").Value = "Sw."
$ws.Range("Q23").Value = 558155
$ws.Range("R23").Value = 7068017
$ws.Range("Z23").Value = $null
$ws.Range("AB23").Value = $null

# Row 24: -> Id 111671226
$ws.Range("A24").Value = 111671226
$ws.Range("B24").Value = 78579
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 2081
$ws.Range("F24").Value = "Skrovellav"
$ws.Range("G24").Value = "Lobaria scrobiculata"
$ws.Range("H24").Value = "(Scop.) DC."
$ws.Range("Q24").Value = 558118
$ws.Range("R24").Value = 7067742
$ws.Range("Z24").Value = $null
$ws.Range("AB24").Value = $null

# Row 25: -> Id 111671294
$ws.Range("A25").Value = 111671294
$ws.Range("B25").Value = 78578
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 6458
$ws.Range("F25").Value = "Lunglav"
$ws.Range("G25").Value = "Lobaria pulmonaria"
$ws.Range("H25").Value = "(L.) Hoffm."
$ws.Range("Q25").Value = 558118
$ws.Range("R25").Value = 7067742
$ws.Range("Z25").Value = $null
$ws.Range("AB25").Value = $null

# Row 26: -> Id 111671190
$ws.Range("A26").Value = 111671190
$ws.Range("B26").Value = 78611
$ws.Range("D26").Value = "LC"
$ws.Range("E26").Value = 6463
$ws.Range("F26").Value = "Bårdlav"
$ws.Range("G26").Value = "Nephroma parile"
$ws.Range("H26").Value = "(Ach.) Ach."
$ws.Range("Q26").Value = 558216
$ws.Range("R26").Value = 7067869
$ws.Range("Z26").Value = $null
$ws.Range("AB26").Value = $null

# Row 27: -> Id 111670497
$ws.Range("A27").Value = 111670497
$ws.Range("B27").Value = 96346
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 620
$ws.Range("F27").Value = "Skogsfru"
$ws.Range("G27").Value = "Epipogium aphyllum"
$ws.Range("H27").Value = "Sw."
$ws.Range("Q27").Value = 558160
$ws.Range("R27").Value = 7068023
$ws.Range("Z27").Value = $null
$ws.Range("AB27").Value = $null

# Row 28: -> Id 111671197
$ws.Range("A28").Value = 111671197
$ws.Range("B28").Value = 78578
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 6458
$ws.Range("F28").Value = "Lunglav"
$ws.Range("G28").Value = "Lobaria pulmonaria"
$ws.Range("H28").Value = "(L.) Hoffm."
$ws.Range("Q28").Value = 558250
$ws.Range("R28").Value = 7067937
$ws.Range("Z28").Value = $null
$ws.Range("AB28").Value = $null

# Row 29: -> Id 111671201
$ws.Range("A29").Value = 111671201
$ws.Range("B29").Value = 78579
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 2081
$ws.Range("F29").Value = "Skrovellav"
$ws.Range("G29").Value = "Lobaria scrobiculata"
$ws.Range("H29").Value = "(Scop.) DC."
$ws.Range("Q29").Value = 558250
$ws.Range("R29").Value = 7067937
$ws.Range("Z29").Value = $null
$ws.Range("AB29").Value = $null

# Row 30: -> Id 111670567
$ws.Range("A30").Value = 111670567
$ws.Range("B30").Value = 96346
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 620
$ws.Range("F30").Value = "Skogsfru"
$ws.Range("G30").Value = "Epipogium aphyllum"
$ws.Range("H30").Value = "Sw."
$ws.Range("Q30").Value = 558130
$ws.Range("R30").Value = 7067959
$ws.Range("Z30").Value = $null
$ws.Range("AB30").Value = $null

# Column L adjustments (Kön / Sex field) - present only for "Skogsfru" (Epipogium aphyllum) rows
# Remove L where no longer applicable
$ws.Range("L19").Value = $null
$ws.Range("L22").Value = $null
$ws.Range("L28").Value = $null

# Add blank L placeholder cells where now applicable (copy blank-cell formatting from existing blank cell)
$ws.Range("I20").Copy($ws.Range("L20"))
$ws.Range("I23").Copy($ws.Range("L23"))
$ws.Range("I27").Copy($ws.Range("L27"))
